$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J45").Value = 0.2858677898194339
$ws.Range("I46").Value = 0.2775335613519331
$ws.Range("H47").Value = 0.2743085116504074
$ws.Range("G48").Value = 0.2534447081011285
$ws.Range("F49").Value = 0.2766837437271186
$ws.Range("E50").Value = 0.2867219094086165
$ws.Range("D51").Value = 0.1751453671933744
$ws.Range("C52").Value = 0.1965658720679752
$ws.Range("B53").Value = 0.4328090033804217
